$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E on this sheet hold plain-text values (prices formatted with
# dotted separators, percentages with padding spaces) stored as inline strings,
# not numbers. Force text format across the data range before writing so
# numeric-looking values (e.g. "596.94") are not auto-converted to real numbers,
# then clear the formatting again so no stray number-format/style is left behind.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '67.476.44'
$ws.Range('E2').Value = '  +0.83%  '
$ws.Range('D3').Value = '3.531.00'
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '596.94'
$ws.Range('E5').Value = '  +0.76%  '
$ws.Range('D6').Value = '173.80'
$ws.Range('E6').Value = '  +2.29%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +2.30%  '
$ws.Range('E9').Value = '  +8.00%  '
$ws.Range('E10').Value = '  +0.70%  '
$ws.Range('D11').Value = '0.437'
$ws.Range('E11').Value = '  -0.13%  '
$ws.Range('D12').Value = '4.142.79'
$ws.Range('E12').Value = '  +0.45%  '
$ws.Range('D14').Value = '28.81'
$ws.Range('E14').Value = '  +2.09%  '
$ws.Range('D15').Value = '0.0000182'
$ws.Range('E15').Value = '  +2.26%  '
$ws.Range('D16').Value = '67.382.52'
$ws.Range('E16').Value = '  +0.73%  '
$ws.Range('D17').Value = '3.543.61'
$ws.Range('E17').Value = '  +0.96%  '
$ws.Range('D18').Value = '6.36'
$ws.Range('E18').Value = '  +1.14%  '
$ws.Range('D19').Value = '14.26'
$ws.Range('E19').Value = '  +1.63%  '
$ws.Range('D20').Value = '398.35'
$ws.Range('E20').Value = '  +1.92%  '
$ws.Range('E21').Value = '  +0.82%  '
$ws.Range('D22').Value = '73.62'
$ws.Range('E22').Value = '  +0.50%  '
$ws.Range('E23').Value = '  +2.40%  '
$ws.Range('E24').Value = '  -0.17%  '
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('E26').Value = '  +1.73%  '
$ws.Range('E27').Value = '  +0.25%  '
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('D29').Value = '6.31'
$ws.Range('E29').Value = '  -1.34%  '
$ws.Range('E30').Value = '  +0.27%  '
$ws.Range('E31').Value = '  +1.16%  '
$ws.Range('E32').Value = '  +2.69%  '
$ws.Range('E33').Value = '  +0.14%  '
$ws.Range('E34').Value = '  +4.48%  '
$ws.Range('D35').Value = '164.26'
$ws.Range('E35').Value = '  +1.67%  '
$ws.Range('E36').Value = '  -1.40%  '
$ws.Range('E37').Value = '  -0.63%  '
$ws.Range('D38').Value = '7.00'
$ws.Range('E38').Value = '  +4.27%  '
$ws.Range('E39').Value = '  +2.43%  '
$ws.Range('D40').Value = '0.0748'
$ws.Range('E40').Value = '  -0.08%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').Value = '26.64'
$ws.Range('E41').Value = '  +0.49%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').Value = '27.41'
$ws.Range('E42').Value = '  +1.81%  '
$ws.Range('E43').Value = '  +3.48%  '
$ws.Range('D44').Value = '2.804.48'
$ws.Range('E44').Value = '  -0.26%  '
$ws.Range('D45').Value = '43.00'
$ws.Range('E45').Value = '  -1.27%  '
$ws.Range('D47').Value = '342.10'
$ws.Range('E47').Value = '  -3.80%  '
$ws.Range('D48').Value = '1.11'
$ws.Range('E48').Value = '  +1.16%  '
$ws.Range('D49').Value = '33.95'
$ws.Range('E49').Value = '  +2.32%  '
$ws.Range('D50').Value = '6.56'
$ws.Range('E50').Value = '  +1.04%  '
$ws.Range('D51').Value = '0.854'
$ws.Range('E51').Value = '  +0.42%  '

$dataRange.ClearFormats()
